# Update the cryptos list (price + volume%) per the GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D look numeric (e.g. "599.84"); assigning them as
# plain Values would let Excel auto-coerce them into Number cells. Prefixing
# with an apostrophe forces a literal-text entry (exactly like a user typing
# '599.84 into the cell) while still storing the text without the quote.
function Set-Price($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
}

Set-Price "D2"  "67.912.70"
$ws.Range("E2").Value  = "  -3.21%  "

Set-Price "D3"  "3.830.75"
$ws.Range("E3").Value  = "  -2.72%  "

$ws.Range("E4").Value  = "  +0.09%  "

Set-Price "D5"  "599.84"
$ws.Range("E5").Value  = "  -1.78%  "

Set-Price "D6"  "167.26"
$ws.Range("E6").Value  = "  -2.39%  "

Set-Price "D7"  "3.830.87"
$ws.Range("E7").Value  = "  -2.64%  "

$ws.Range("E8").Value  = "  +0.16%  "

$ws.Range("E9").Value  = "  -1.82%  "

$ws.Range("E10").Value = "  -4.04%  "

$ws.Range("E11").Value = "  +0.48%  "

Set-Price "D12" "0.458"
$ws.Range("E12").Value = "  -2.67%  "

Set-Price "D13" "0.0000259"
$ws.Range("E13").Value = "  +0.75%  "

Set-Price "D14" "37.08"
$ws.Range("E14").Value = "  -4.07%  "

Set-Price "D15" "4.478.61"
$ws.Range("E15").Value = "  -2.55%  "

Set-Price "D16" "3.843.43"
$ws.Range("E16").Value = "  -2.52%  "

Set-Price "D17" "68.108.34"
$ws.Range("E17").Value = "  -2.92%  "

Set-Price "D18" "18.26"
$ws.Range("E18").Value = "  -1.77%  "

Set-Price "D19" "7.39"
$ws.Range("E19").Value = "  -3.88%  "

$ws.Range("E20").Value = "  -0.82%  "

Set-Price "D21" "11.09"
$ws.Range("E21").Value = "  -0.43%  "

Set-Price "D22" "466.02"
$ws.Range("E22").Value = "  -6.30%  "

Set-Price "D23" "0.734"
$ws.Range("E23").Value = "  -2.04%  "

$ws.Range("E24").Value = "  -3.97%  "

Set-Price "D25" "82.88"
$ws.Range("E25").Value = "  -3.81%  "

$ws.Range("E26").Value = "  -3.05%  "

Set-Price "D27" "12.12"
$ws.Range("E27").Value = "  -2.29%  "

Set-Price "D28" "10.04"
$ws.Range("E28").Value = "  -1.69%  "

Set-Price "D29" "1.00"
$ws.Range("E29").Value = "  +0.06%  "

Set-Price "D31" "3.984.54"
$ws.Range("E31").Value = "  -2.49%  "

Set-Price "D32" "7.64"
$ws.Range("E32").Value = "  -3.34%  "

Set-Price "D33" "2.31"
$ws.Range("E33").Value = "  -6.22%  "

Set-Price "D34" "31.27"
$ws.Range("E34").Value = "  -3.42%  "

$ws.Range("E35").Value = "  -0.85%  "

Set-Price "D36" "3.797.11"
$ws.Range("E36").Value = "  -2.63%  "

$ws.Range("E37").Value = "  -3.53%  "

Set-Price "D38" "3.65"
$ws.Range("E38").Value = "  +10.66%  "

$ws.Range("E39").Value = "  -1.29%  "

$ws.Range("E40").Value = "  -2.88%  "

Set-Price "D41" "5.92"
$ws.Range("E41").Value = "  -4.17%  "

$ws.Range("E42").Value = "  +0.14%  "

Set-Price "D43" "0.314"
$ws.Range("E43").Value = "  -4.98%  "

$ws.Range("E44").Value = "  -6.53%  "

Set-Price "D45" "420.43"
$ws.Range("E45").Value = "  -4.15%  "

Set-Price "D46" "8.72"
$ws.Range("E46").Value = "  +0.33%  "

$ws.Range("E48").Value = "  +5.89%  "

Set-Price "D49" "46.89"
$ws.Range("E49").Value = "  -3.00%  "

Set-Price "D50" "26.39"
$ws.Range("E50").Value = "  +3.89%  "

Set-Price "D51" "142.12"
$ws.Range("E51").Value = "  -0.99%  "
